# Fix: blank first row in Excel was being dropped when exporting to PDF.
# Put a (visually unobtrusive, but non-empty) value in A1 so the row is no
# longer entirely blank, and format it to match the header band look:
# bold, 12pt, red, Times New Roman - keeping the cell's existing green
# fill / vertical-center + wrap alignment untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A1")
$cell.Value = "A1"

$cell.Font.Name = "Times New Roman"
$cell.Font.Size = 12
$cell.Font.Bold = $true
$cell.Font.Color = 255   # VBA BGR long for RGB(255,0,0) -> OOXML rgb="FFFF0000"
